$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.720.15'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.622.62'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.95'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.14%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.63'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').Value = '3.092.42'
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('D15').Value = '63.535.69'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').Value = '2.633.53'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.26'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +7.00%  '
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '348.30'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.97%  '
$ws.Range('E21').Value = '  -1.02%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  +2.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.26'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +11.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.16'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.67'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '563.67'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.21'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.54%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').Value = '0.0₃0841'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.21'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '168.71'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '169.21'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.88'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.90'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0595'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.29'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.98'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.78%  '
$ws.Range('E50').Value = '  +0.81%  '
$ws.Range('E51').Value = '  +1.99%  '
